$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 12.12
$ws.Range("E2").Value = 19.998
$ws.Range("F2").Value = -39.996
$ws.Range("G2").Value = -19.998
$ws.Range("D3").Value = 13.41
$ws.Range("E3").Value = 41.68550140394704
$ws.Range("F3").Value = -81.67053737332458
$ws.Range("G3").Value = -39.98503596937753
$ws.Range("D4").Value = 14.76
$ws.Range("E4").Value = 65.12337382404914
$ws.Range("F4").Value = -125.1113581840309
$ws.Range("G4").Value = -59.98798435998179
$ws.Range("D5").Value = 16.14
$ws.Range("E5").Value = 90.32314913328499
$ws.Range("F5").Value = -170.3045494879536
$ws.Range("G5").Value = -79.9814003546686
$ws.Range("B6").Value = "jog(4分59秒~4分)"
$ws.Range("C6").Value = 1.5
$ws.Range("D6").Value = 15.16
$ws.Range("E6").Value = 113.3521057671903
$ws.Range("F6").Value = -209.3490956286434
$ws.Range("G6").Value = -95.99698986145306
$ws.Range("E7").Value = 145.8609520846057
$ws.Range("F7").Value = -265.847541270648
$ws.Range("G7").Value = -119.9865891860422
$ws.Range("B8").Value = "jog(~5分)"
$ws.Range("D8").Value = 6.7
$ws.Range("E8").Value = 153.7103473411248
$ws.Range("F8").Value = -270.8122318121725
$ws.Range("G8").Value = -117.1018844710477
$ws.Range("E9").Value = 150.3322354371725
$ws.Range("F9").Value = -253.346734492289
$ws.Range("G9").Value = -103.0144990551165
$ws.Range("E10").Value = 147.0283647292947
$ws.Range("F10").Value = -237.0076397524871
$ws.Range("G10").Value = -89.97927502319232
$ws.Range("B11").Value = "jog(4分59秒~4分)"
$ws.Range("C11").Value = 1.5
$ws.Range("D11").Value = 13.37
$ws.Range("E11").Value = 165.8581036092584
$ws.Range("F11").Value = -265.844302494309
$ws.Range("G11").Value = -99.98619888505061
$ws.Range("D12").Value = 14.24
$ws.Range("E12").Value = 185.7090189167754
$ws.Range("F12").Value = -295.6912019142844
$ws.Range("G12").Value = -109.982182997509
$ws.Range("D13").Value = 15.15
$ws.Range("E13").Value = 206.6256681272816
$ws.Range("F13").Value = -326.6171848032075
$ws.Range("G13").Value = -119.9915166759258
$ws.Range("D14").Value = 17.09
$ws.Range("E14").Value = 220.8836294708938
$ws.Range("F14").Value = -343.1506578147633
$ws.Range("G14").Value = -122.2670283438696
$ws.Range("E15").Value = 251.0292417799477
$ws.Range("F15").Value = -391.0198373039056
$ws.Range("G15").Value = -139.9905955239579
